$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 189, shifting all data at/after
# row 189 down by two rows (so old row189/190 become row191/192, etc, and
# the two rows that fell off the bottom become new rows 268/269).
$ws.Rows("189:190").Insert()

# Fill in the two newly inserted rows with their data.
# Row 189 (new record)
$ws.Range("A189").Value = 7
$ws.Range("B189").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C189").Value = "Ñuble"
$ws.Range("D189").Value = 44460
$ws.Range("E189").Value = 16
$ws.Range("F189").Value = "Fruta"
$ws.Range("G189").Value = 100106
$ws.Range("H189").Value = "Oleaginosos"
$ws.Range("I189").Value = 100106002
$ws.Range("J189").Value = "Palta"
$ws.Range("K189").Value = "Hass"
$ws.Range("L189").Value = "1a nueva(o)"
$ws.Range("M189").Value = 120
$ws.Range("N189").Value = 2700
$ws.Range("O189").Value = 2800
$ws.Range("P189").Value = 2750
$ws.Range("Q189").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R189").Value = "Provincia de Quillota"
$ws.Range("S189").Value = 2750
$ws.Range("T189").Value = 1

# Row 190 (new record)
$ws.Range("A190").Value = 7
$ws.Range("B190").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C190").Value = "Ñuble"
$ws.Range("D190").Value = 44460
$ws.Range("E190").Value = 16
$ws.Range("F190").Value = "Fruta"
$ws.Range("G190").Value = 100106
$ws.Range("H190").Value = "Oleaginosos"
$ws.Range("I190").Value = 100106002
$ws.Range("J190").Value = "Palta"
$ws.Range("K190").Value = "Hass"
$ws.Range("L190").Value = "Primera"
$ws.Range("M190").Value = 400
$ws.Range("N190").Value = 25000
$ws.Range("O190").Value = 26000
$ws.Range("P190").Value = 25500
$ws.Range("Q190").Value = "$/bandeja 10 kilos"
$ws.Range("R190").Value = "Perú"
$ws.Range("S190").Value = 2550
$ws.Range("T190").Value = 10
